# Applies the cell-value updates described in the commit diff
# ("Update gh-pages to output generated at 456a3b4") across all four
# worksheets of the workbook: 展览, 演出, 本地生活, 全部类型.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value2 = 9026
$ws.Range("G4").Value2 = 60
$ws.Range("G5").Value2 = "不可售"
$ws.Range("F6").Value2 = 743
$ws.Range("F8").Value2 = 83
$ws.Range("F10").Value2 = 920
$ws.Range("F11").Value2 = 4018
$ws.Range("F12").Value2 = 319
$ws.Range("F13").Value2 = 198
$ws.Range("G13").Value2 = 59
$ws.Range("F14").Value2 = 816
$ws.Range("F19").Value2 = 24
$ws.Range("F20").Value2 = 1456
$ws.Range("F21").Value2 = 1373
$ws.Range("F22").Value2 = 534
$ws.Range("F25").Value2 = 187
$ws.Range("F26").Value2 = 391
$ws.Range("F27").Value2 = 78
$ws.Range("F31").Value2 = 780
$ws.Range("F32").Value2 = 80
$ws.Range("F33").Value2 = 63
$ws.Range("F34").Value2 = 118
$ws.Range("F35").Value2 = 8
$ws.Range("F38").Value2 = 224
$ws.Range("F39").Value2 = 208
$ws.Range("F40").Value2 = 431
$ws.Range("F42").Value2 = 33

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value2 = 85
$ws.Range("F6").Value2 = 63

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value2 = 223

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value2 = 223
$ws.Range("F3").Value2 = 9026
$ws.Range("G4").Value2 = 60
$ws.Range("G5").Value2 = "不可售"
$ws.Range("F6").Value2 = 743
$ws.Range("F8").Value2 = 83
$ws.Range("F10").Value2 = 920
$ws.Range("F12").Value2 = 4018
$ws.Range("F13").Value2 = 319
$ws.Range("F14").Value2 = 198
$ws.Range("G14").Value2 = 59
$ws.Range("F16").Value2 = 85
$ws.Range("F17").Value2 = 816
$ws.Range("F20").Value2 = 63
$ws.Range("F25").Value2 = 24
$ws.Range("F26").Value2 = 1456
$ws.Range("F27").Value2 = 1373
$ws.Range("F28").Value2 = 534
$ws.Range("F31").Value2 = 187
$ws.Range("F33").Value2 = 391
$ws.Range("F34").Value2 = 78
$ws.Range("F37").Value2 = 780
$ws.Range("F38").Value2 = 80
$ws.Range("F39").Value2 = 63
$ws.Range("F40").Value2 = 118
$ws.Range("F41").Value2 = 8
$ws.Range("F44").Value2 = 208
$ws.Range("F45").Value2 = 431
$ws.Range("F47").Value2 = 33
